$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. Update "Valor Mora" total (E11): 300000 -> 400000
$ws.Range("E11").Value = 400000

# 2. Update "Cant. Periodos" (F13): 3 -> 4
$ws.Range("F13").Value = 4

# 3. Insert a new row for period "2509", copying the format of the last
#    data row (row 18) so the table keeps its borders/shading, then push
#    the existing last row's format back to "middle" style.
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B19:J19").Insert(-4121) | Out-Null  # xlShiftDown, using copied formatting
$excel.CutCopyMode = 0

# Fill in the new row's data (period 2509, same worker/amounts as others)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45555876"
$ws.Range("D19").Value = "VIVIANA PAOLA MESTRA PADILLA"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 100000
$ws.Range("G19").Value = 2500000

# Row 18 should now look like the other "middle" rows (15:18) rather than
# the table's final bordered row, so copy that formatting down from row 17.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
